$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.380.70'
$ws.Range("E2").Value = '  -3.08%  '
$ws.Range("D3").Value = '1.775.40'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4233'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.56%  '
$ws.Range("E8").Value = '  +1.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07125'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8357'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.36'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("D12").Value = '1.800.28'
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.442'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.231'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06867'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008670'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.78%  '
$ws.Range("D21").Value = '26.389.04'
$ws.Range("E21").Value = '  -2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.065'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.16%  '
$ws.Range("D24").Value = '1.995.89'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.791'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -9.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.061'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.818'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08845'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7251'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.74%  '
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.313'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.73%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.732'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.085'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05130'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01881'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1609'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4902'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.607'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.360'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.955'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.78%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.19'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.635'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06172'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4428'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.715'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.08%  '
